$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column A, shifting existing data right (A:D -> D:G)
$ws.Range("A1:C1").EntireColumn.Insert()

# Fill the newly inserted header cells with the new values
$ws.Range("A1").Value = "test_22@gmail.com"
$ws.Range("B1").Value = "test_5@gmail.com"
$ws.Range("C1").Value = "test@gmail.com"

# Match the header formatting (bold/border/center) used by the other header cells
$ws.Range("D1").Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)
